$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh - price (D) and 1h volume (E) updates,
# plus a few coins that swapped rank position (B/C/D/E all change).

# Row 2
$ws.Cells.Item(2, 4).Value = '59.261.15'
$ws.Cells.Item(2, 5).Value = '  -2.49%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.577.99'
$ws.Cells.Item(3, 5).Value = '  -2.34%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '555.61'
$ws.Cells.Item(5, 5).Value = '  -2.06%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '142.13'
$ws.Cells.Item(6, 5).Value = '  -3.29%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.11%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.597'
$ws.Cells.Item(8, 5).Value = '  -2.25%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.580.56'
$ws.Cells.Item(9, 5).Value = '  -3.08%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '6.64'
$ws.Cells.Item(10, 5).Value = '  -3.18%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.104'
$ws.Cells.Item(11, 5).Value = '  -1.15%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +11.63%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.356'
$ws.Cells.Item(13, 5).Value = '  +3.76%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.034.23'
$ws.Cells.Item(14, 5).Value = '  -2.80%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '59.272.04'
$ws.Cells.Item(15, 5).Value = '  -2.35%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '23.12'
$ws.Cells.Item(16, 5).Value = '  +5.06%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.0000137'
$ws.Cells.Item(17, 5).Value = '  -0.58%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.582.91'
$ws.Cells.Item(18, 5).Value = '  -2.74%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.56'
$ws.Cells.Item(19, 5).Value = '  +0.28%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '337.40'
$ws.Cells.Item(20, 5).Value = '  -1.74%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.37'
$ws.Cells.Item(21, 5).Value = '  -1.02%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.48'
$ws.Cells.Item(22, 5).Value = '  +1.74%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.09%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '63.29'
$ws.Cells.Item(24, 5).Value = '  -5.48%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.471'
$ws.Cells.Item(25, 5).Value = '  +6.47%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  +0.27%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.161'
$ws.Cells.Item(27, 5).Value = '  -2.40%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.47'
$ws.Cells.Item(28, 5).Value = '  +0.80%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '0.0₃0774'
$ws.Cells.Item(29, 5).Value = '  -4.29%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  +0.06%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '6.18'
$ws.Cells.Item(31, 5).Value = '  -1.32%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.67'
$ws.Cells.Item(32, 5).Value = '  -2.89%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '158.47'
$ws.Cells.Item(33, 5).Value = '  +0.27%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '19.10'
$ws.Cells.Item(34, 5).Value = '  -0.80%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.05'
$ws.Cells.Item(35, 5).Value = '  -1.31%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.25%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.900'
$ws.Cells.Item(37, 5).Value = '  +0.15%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'OKB'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '37.26'
$ws.Cells.Item(38, 5).Value = '  -0.61%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'SuiNetwork'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.858'
$ws.Cells.Item(39, 5).Value = '  -5.67%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.48'
$ws.Cells.Item(40, 5).Value = '  -2.36%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Filecoin'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.67'
$ws.Cells.Item(41, 5).Value = '  +0.40%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '292.55'
$ws.Cells.Item(42, 5).Value = '  -3.45%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '136.60'
$ws.Cells.Item(43, 5).Value = '  +5.98%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.998'
$ws.Cells.Item(44, 5).Value = '  +0.09%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0976'
$ws.Cells.Item(45, 5).Value = '  -1.23%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.591'
$ws.Cells.Item(46, 5).Value = '  -2.39%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.66'
$ws.Cells.Item(47, 5).Value = '  -0.43%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0532'
$ws.Cells.Item(48, 5).Value = '  -2.70%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0235'
$ws.Cells.Item(49, 5).Value = '  -1.31%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '18.70'
$ws.Cells.Item(50, 5).Value = '  -0.70%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '1.957.53'
$ws.Cells.Item(51, 5).Value = '  -0.37%  '

